# Apply the "page updated for advertisement" schedule edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the coffee break topic (shared across the three day coffee-break rows)
#    and the "Advanced and combining plots" topic on day 1 -> "Advanced ggplot".
$ws.Range("E4").Value = "Break"
$ws.Range("E5").Value = "Advanced ggplot"
$ws.Range("E13").Value = "Break"
$ws.Range("E23").Value = "Break"

# 2. Remove the standalone "Break" row on day 2 (row 20, 16:60-18:30) -
#    the dinner / Restaurang Valvet row moves up to take its place.
$ws.Rows(20).Delete()

# Leave the active selection on E5, matching where the topic text was edited.
$ws.Range("E5").Select()
